# The workbook holds a weekly price log (Orégano, Mercado Mayorista Lo
# Valledor de Santiago). This commit ("Fruta / hortaliza, semanal") adds one
# new weekly record at the top of the data block (row 13, right after the
# last "recent" summary rows 2-12), pushing all existing records from row 13
# down to row 14, and growing the sheet from A1:R116 to A1:R117.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 13 - shifts rows 13..116 down to 14..117
# and Excel extends the used range / dimension automatically.
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new weekly record.
$ws.Cells.Item(13, 1).Value  = 6
$ws.Cells.Item(13, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(13, 3).Value  = "Metropolitana"
$ws.Cells.Item(13, 4).Value  = 44490
$ws.Cells.Item(13, 5).Value  = 13
$ws.Cells.Item(13, 6).Value  = 100112029
$ws.Cells.Item(13, 7).Value  = "Orégano"
$ws.Cells.Item(13, 8).Value  = "Sin especificar"
$ws.Cells.Item(13, 9).Value  = "Primera"
$ws.Cells.Item(13, 10).Value = 34
$ws.Cells.Item(13, 11).Value = 8500
$ws.Cells.Item(13, 12).Value = 9000
$ws.Cells.Item(13, 13).Value = 8735
$ws.Cells.Item(13, 14).Value = "`$/docena de atados"
$ws.Cells.Item(13, 15).Value = "Región Metropolitana"
$ws.Cells.Item(13, 16).Value = 2912
$ws.Cells.Item(13, 17).Value = 3
$ws.Cells.Item(13, 18).Value = "Hortaliza"
